# Update the "Neg_Change" worksheet with the new filtered data rows
$wb = $excel.ActiveWorkbook

$negRows = @(
    @("BAJAJFINSV", 2009.5, 2039, 1995, 2039, 2396251, 4943745, -0.5152963997940833, "BAJAJFINSV"),
    @("BAJFINANCE", 1011, 1018.4, 1008.8, 1011, 11131278, 23762890, -0.5315688453719224, "BAJFINANCE"),
    @("ZYDUSLIFE", 940.1, 957.25, 940, 951.1, 791415, 1819775, -0.5651028286463986, "ZYDUSLIFE"),
    @("SYNGENE", 649.3, 667.4, 640.65, 666, 1155741, 2473758, -0.5327994896833077, "SYNGENE"),
    @("JUBLFOOD", 581.05, 586.65, 578, 584.65, 1226844, 2599968, -0.5281311154598826, "JUBLFOOD"),
    @("BHEL", 285.5, 290.9, 282.75, 286.7, 15184051, 37222053, -0.5920684170752215, "BHEL"),
    @("NATIONALUM", 265.5, 269.25, 264.2, 266.1, 17311954, 36187656, -0.5216060968414202, "NATIONALUM"),
    @("KALYANKJIL", 510.25, 515, 505.7, 508.7, 1756494, 4125402, -0.5742247664591232, "KALYANKJIL"),
    @("PHOENIXLTD", 1738, 1744.4, 1718.8, 1731.9, 238006, 511582, -0.534764710251729, "PHOENIXLTD"),
    @("ABCAPITAL", 329.35, 333.1, 327.9, 327.9, 2257142, 4844131, -0.5340460445846738, "ABCAPITAL"),
    @("IIFL", 543.9, 544.8, 534.2, 540, 566485, 1213563, -0.5332051158448304, "IIFL")
)

$posRows = @(
    @("TCS", 3060, 3136.4, 3057.1, 3125.1, 4460315, 2806868, 0.5890718765542234, "TCS"),
    @("WIPRO", 242.75, 246.3, 242.1, 245.39, 13152546, 8831094, 0.4893450347148383, "WIPRO"),
    @("ULTRACEMCO", 11861, 11949, 11831, 11920, 184542, 125847, 0.4663996757968009, "ULTRACEMCO"),
    @("HINDALCO", 798, 809, 791.75, 793.95, 7589035, 5411465, 0.4023993502683654, "HINDALCO"),
    @("POWERGRID", 269, 270.25, 266.5, 266.85, 16444701, 11154806, 0.4742256387067601, "POWERGRID"),
    @("BOSCHLTD", 36610, 37345, 35750, 37160, 53670, 36836, 0.4569985883374959, "BOSCHLTD"),
    @("AMBUJACEM", 560, 570.05, 557.15, 563, 1748539, 1177397, 0.4850887168898851, "AMBUJACEM"),
    @("PIDILITIND", 1477, 1487, 1469.7, 1481, 417016, 292665, 0.4248919412980712, "PIDILITIND"),
    @("DLF", 767.4, 771.15, 760, 762, 2364484, 1497577, 0.5788730729705384, "DLF"),
    @("DMART", 4112.6, 4112.6, 4047.3, 4054, 381769, 272277, 0.4021345908762033, "DMART"),
    @("ABB", 4994, 5000, 4911, 4933, 322415, 202095, 0.5953635666394518, "ABB"),
    @("VBL", 470.6, 470.6, 458.6, 459.2, 6996274, 4669423, 0.4983166014301981, "VBL"),
    @("HEROMOTOCO", 5440, 5557, 5421.5, 5530, 534031, 348232, 0.5335494727652829, "HEROMOTOCO"),
    @("AUROPHARMA", 1170, 1195.5, 1160.9, 1190.3, 2367201, 1481205, 0.5981589314105745, "AUROPHARMA"),
    @("OFSS", 8300, 8430, 8271.5, 8393, 104111, 74067, 0.405632737926472, "OFSS"),
    @("CONCOR", 530, 538.5, 524.7, 531, 4975992, 3241835, 0.5349306796922113, "CONCOR"),
    @("DIXON", 15155, 15369, 15054, 15300, 260843, 183854, 0.4187507478760321, "DIXON"),
    @("IREDA", 149.76, 153.48, 149.55, 150.96, 6667190, 4634146, 0.4387095270628072, "IREDA"),
    @("KEI", 4074.9, 4098.2, 4062, 4073.6, 133502, 90590, 0.4736946682856827, "KEI"),
    @("CUMMINSIND", 4420, 4423.9, 4370, 4385, 617858, 395014, 0.5641420304090488, "CUMMINSIND"),
    @("ASHOKLEY", 147.47, 147.6, 141.74, 142.22, 24925913, 17177183, 0.451105981696766, "ASHOKLEY"),
    @("PGEL", 534, 536.05, 519.55, 527, 1917942, 1294251, 0.4818933885312818, "PGEL")
)

$wsNeg = $wb.Worksheets.Item("Neg_Change")
for ($i = 0; $i -lt $negRows.Count; $i++) {
    $r = $i + 2
    $row = $negRows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $wsNeg.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

$wsPos = $wb.Worksheets.Item("Pos_Change")
for ($i = 0; $i -lt $posRows.Count; $i++) {
    $r = $i + 2
    $row = $posRows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $wsPos.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
